$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and 1h volume change (E) columns with refreshed crypto data.
# Force text storage (matching the source inlineStr cells) for numeric-looking
# price strings so Excel does not silently coerce them into numbers.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.141.23"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -0.53%  "
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.631.14"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -1.16%  "
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.10%  "
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "216.23"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -1.10%  "
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.516"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  +0.91%  "
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +0.12%  "
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -1.47%  "
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -0.82%  "
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -0.60%  "
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0851"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  +0.34%  "
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.860.34"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -1.06%  "
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.637.38"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -0.57%  "
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -0.63%  "
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.539"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  +0.16%  "
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "65.31"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -3.58%  "
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "27.125.84"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -0.43%  "
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0731"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -1.64%  "
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "213.77"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -3.35%  "
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  +0.88%  "
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.37"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -1.81%  "
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -0.85%  "
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "9.09"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -1.77%  "
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "147.05"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -1.04%  "
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  +0.14%  "
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "7.30"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -1.63%  "
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -1.08%  "
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "15.55"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -1.82%  "
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -0.54%  "
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -0.90%  "
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.37"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +0.22%  "
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -1.43%  "
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.317.77"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +3.14%  "
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -1.77%  "
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -0.26%  "
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -1.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.539"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -0.92%  "
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -0.61%  "
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +0.01%  "
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  +2.18%  "
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -1.02%  "
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.26"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -2.15%  "
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.768.91"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -1.19%  "
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "62.55"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -0.45%  "
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "90.83"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -1.88%  "
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -0.02%  "
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -0.92%  "
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.807"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  +19.34%  "
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0513"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -0.22%  "
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "7.57"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -1.77%  "
$cell.Style = "Normal"
